# Scheduled-runner market data refresh: updates currentAveragePrice* /
# LevePrice* / LeveProfit* columns (H:N) for the affected leve rows on
# each job-sheet, mirroring the latest Universalis price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 9.6
$ws.Range("I11").Value = 9.6
$ws.Range("K11").Value = 9.6
$ws.Range("M11").Value = 130.4

$ws.Range("H80").Value = 844.2
$ws.Range("I80").Value = 531.8889
$ws.Range("J80").Value = 1099.7273
$ws.Range("K80").Value = 1595.6667
$ws.Range("L80").Value = 3299.1819
$ws.Range("M80").Value = -597.6667000000002
$ws.Range("N80").Value = -5295.1819

$ws.Range("H83").Value = 844.2
$ws.Range("I83").Value = 531.8889
$ws.Range("J83").Value = 1099.7273
$ws.Range("K83").Value = 4787.0001
$ws.Range("L83").Value = 9897.545700000001
$ws.Range("M83").Value = 204.9998999999998
$ws.Range("N83").Value = -19881.5457

$ws.Range("H86").Value = 2859970.8
$ws.Range("I86").Value = 3112.7144
$ws.Range("J86").Value = 9525973
$ws.Range("K86").Value = 3112.7144
$ws.Range("L86").Value = 9525973
$ws.Range("M86").Value = -1989.7144
$ws.Range("N86").Value = -9528219

$ws.Range("H89").Value = 2859970.8
$ws.Range("I89").Value = 3112.7144
$ws.Range("J89").Value = 9525973
$ws.Range("K89").Value = 15563.572
$ws.Range("L89").Value = 47629865
$ws.Range("M89").Value = -9947.572
$ws.Range("N89").Value = -47641097

$ws.Range("H111").Value = 2056.5
$ws.Range("J111").Value = 3092
$ws.Range("L111").Value = 9276
$ws.Range("N111").Value = -15410

$ws.Range("H113").Value = 4805.1953
$ws.Range("I113").Value = 6275
$ws.Range("J113").Value = 4646.2974
$ws.Range("K113").Value = 6275
$ws.Range("L113").Value = 4646.2974
$ws.Range("M113").Value = -3021
$ws.Range("N113").Value = -11154.2974

$ws.Range("H137").Value = 1880.0513
$ws.Range("J137").Value = 1543.625
$ws.Range("L137").Value = 4630.875
$ws.Range("N137").Value = -9730.875

$ws.Range("H138").Value = 2262.5386
$ws.Range("J138").Value = 2660.0715
$ws.Range("L138").Value = 7980.2145
$ws.Range("N138").Value = -18260.2145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3893.8044
$ws.Range("I32").Value = 2454.0588
$ws.Range("K32").Value = 2454.0588
$ws.Range("M32").Value = -2167.0588

$ws.Range("H45").Value = 2125.3333
$ws.Range("I45").Value = 1981.2
$ws.Range("K45").Value = 1981.2
$ws.Range("M45").Value = -1604.2

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

$ws.Range("H88").Value = 3270
$ws.Range("I88").Value = 2721.8
$ws.Range("J88").Value = 3498.4167
$ws.Range("K88").Value = 2721.8
$ws.Range("L88").Value = 3498.4167
$ws.Range("M88").Value = -2315.8
$ws.Range("N88").Value = -4310.4167

$ws.Range("H91").Value = 3270
$ws.Range("I91").Value = 2721.8
$ws.Range("J91").Value = 3498.4167
$ws.Range("K91").Value = 2721.8
$ws.Range("L91").Value = 3498.4167
$ws.Range("M91").Value = -1317.8
$ws.Range("N91").Value = -6306.4167

$ws.Range("H132").Value = 12438
$ws.Range("I132").Value = 6491.893
$ws.Range("J132").Value = 67935
$ws.Range("K132").Value = 19475.679
$ws.Range("L132").Value = 203805
$ws.Range("M132").Value = -16945.679
$ws.Range("N132").Value = -208865

$ws.Range("H135").Value = 59950
$ws.Range("J135").Value = 59950
$ws.Range("L135").Value = 59950
$ws.Range("N135").Value = -70090

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 71433710
$ws.Range("I99").Value = 83339130
$ws.Range("K99").Value = 83339130
$ws.Range("M99").Value = -83337632

$ws.Range("H107").Value = 9719.52
$ws.Range("I107").Value = 9135.227999999999
$ws.Range("K107").Value = 9135.227999999999
$ws.Range("M107").Value = -7215.227999999999

$ws.Range("H134").Value = 2891.103
$ws.Range("I134").Value = 2608.389
$ws.Range("J134").Value = 3981.5715
$ws.Range("K134").Value = 7825.167
$ws.Range("L134").Value = 11944.7145
$ws.Range("M134").Value = -5290.167
$ws.Range("N134").Value = -17014.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1542.3214
$ws.Range("I31").Value = 1145.9524
$ws.Range("J31").Value = 2731.4285
$ws.Range("K31").Value = 1145.9524
$ws.Range("L31").Value = 2731.4285
$ws.Range("M31").Value = -850.9523999999999
$ws.Range("N31").Value = -3321.4285

$ws.Range("H34").Value = 1542.3214
$ws.Range("I34").Value = 1145.9524
$ws.Range("J34").Value = 2731.4285
$ws.Range("K34").Value = 1145.9524
$ws.Range("L34").Value = 2731.4285
$ws.Range("M34").Value = -943.9523999999999
$ws.Range("N34").Value = -3135.4285

$ws.Range("H99").Value = 6783
$ws.Range("I99").Value = 3912
$ws.Range("J99").Value = 7500.75
$ws.Range("K99").Value = 3912
$ws.Range("L99").Value = 7500.75
$ws.Range("M99").Value = -2414
$ws.Range("N99").Value = -10496.75

$ws.Range("H126").Value = 6783
$ws.Range("I126").Value = 3912
$ws.Range("J126").Value = 7500.75
$ws.Range("K126").Value = 11736
$ws.Range("L126").Value = 22502.25
$ws.Range("M126").Value = -9266
$ws.Range("N126").Value = -27442.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 106350540
$ws.Range("J46").Value = 119049280
$ws.Range("L46").Value = 357147840
$ws.Range("N46").Value = -357148022

$ws.Range("H122").Value = 1169.05
$ws.Range("I122").Value = 2821.5
$ws.Range("J122").Value = 755.9375
$ws.Range("K122").Value = 25393.5
$ws.Range("L122").Value = 6803.4375
$ws.Range("M122").Value = -22943.5
$ws.Range("N122").Value = -11703.4375

$ws.Range("H132").Value = 966
$ws.Range("I132").Value = 615
$ws.Range("K132").Value = 5535
$ws.Range("M132").Value = -3005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2488.4849
$ws.Range("I122").Value = 2012.6522
$ws.Range("K122").Value = 6037.9566
$ws.Range("M122").Value = -3587.9566

$ws.Range("H126").Value = 5870.857
$ws.Range("I126").Value = 5019.4
$ws.Range("J126").Value = 7999.5
$ws.Range("K126").Value = 15058.2
$ws.Range("L126").Value = 23998.5
$ws.Range("M126").Value = -12588.2
$ws.Range("N126").Value = -28938.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4468.0435
$ws.Range("I7").Value = 3395.4285
$ws.Range("J7").Value = 4937.3125
$ws.Range("K7").Value = 3395.4285
$ws.Range("L7").Value = 4937.3125
$ws.Range("M7").Value = -3283.4285
$ws.Range("N7").Value = -5161.3125

$ws.Range("H40").Value = 5227.375
$ws.Range("I40").Value = 5028.647
$ws.Range("J40").Value = 5710
$ws.Range("K40").Value = 5028.647
$ws.Range("L40").Value = 5710
$ws.Range("M40").Value = -4892.647
$ws.Range("N40").Value = -5982

$ws.Range("H82").Value = 3372.8333
$ws.Range("I82").Value = 1674.7858
$ws.Range("J82").Value = 9316
$ws.Range("K82").Value = 1674.7858
$ws.Range("L82").Value = 9316
$ws.Range("M82").Value = -1313.7858
$ws.Range("N82").Value = -10038

$ws.Range("H85").Value = 3372.8333
$ws.Range("I85").Value = 1674.7858
$ws.Range("J85").Value = 9316
$ws.Range("K85").Value = 1674.7858
$ws.Range("L85").Value = 9316
$ws.Range("M85").Value = -426.7858000000001
$ws.Range("N85").Value = -11812

$ws.Range("H126").Value = 4468.0435
$ws.Range("I126").Value = 3395.4285
$ws.Range("J126").Value = 4937.3125
$ws.Range("K126").Value = 10186.2855
$ws.Range("L126").Value = 14811.9375
$ws.Range("M126").Value = -7716.2855
$ws.Range("N126").Value = -19751.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2462
$ws.Range("I126").Value = 2448
$ws.Range("K126").Value = 7344
$ws.Range("M126").Value = -4874
